$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 84.72727
$ws.Range("I9").Value = 92.2
$ws.Range("J9").Value = 10
$ws.Range("K9").Value = 92.2
$ws.Range("L9").Value = 10
$ws.Range("M9").Value = 76.8
$ws.Range("N9").Value = -348
$ws.Range("H93").Value = 601000000
$ws.Range("J93").Value = 601000000
$ws.Range("L93").Value = 601000000
$ws.Range("N93").Value = -601004992
$ws.Range("H100").Value = 60148.28
$ws.Range("I100").Value = 63064.312
$ws.Range("J100").Value = 54964.223
$ws.Range("K100").Value = 63064.312
$ws.Range("L100").Value = 54964.223
$ws.Range("M100").Value = -62523.312
$ws.Range("N100").Value = -56046.223
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H113").Value = 10378.223
$ws.Range("I113").Value = 9600.799999999999
$ws.Range("K113").Value = 9600.799999999999
$ws.Range("M113").Value = -6346.799999999999
$ws.Range("H131").Value = 2122.2666
$ws.Range("I131").Value = 1659.4546
$ws.Range("K131").Value = 4978.3638
$ws.Range("M131").Value = 61.63619999999992
$ws.Range("H137").Value = 11147.318
$ws.Range("I137").Value = 12902.833
$ws.Range("K137").Value = 38708.499
$ws.Range("M137").Value = -36158.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9785.02
$ws.Range("I32").Value = 9834.83
$ws.Range("J32").Value = 9199.75
$ws.Range("K32").Value = 9834.83
$ws.Range("L32").Value = 9199.75
$ws.Range("M32").Value = -9547.83
$ws.Range("N32").Value = -9773.75
$ws.Range("H61").Value = 8193.379000000001
$ws.Range("I61").Value = 9178.913
$ws.Range("J61").Value = 4415.5
$ws.Range("K61").Value = 9178.913
$ws.Range("L61").Value = 4415.5
$ws.Range("M61").Value = -8966.913
$ws.Range("N61").Value = -4839.5
$ws.Range("H74").Value = 4374.343
$ws.Range("I74").Value = 5962.7896
$ws.Range("K74").Value = 5962.7896
$ws.Range("M74").Value = -5088.7896
$ws.Range("H77").Value = 4374.343
$ws.Range("I77").Value = 5962.7896
$ws.Range("K77").Value = 29813.948
$ws.Range("M77").Value = -25445.948
$ws.Range("H132").Value = 2235.5686
$ws.Range("I132").Value = 1699.4572
$ws.Range("J132").Value = 3408.3125
$ws.Range("K132").Value = 5098.3716
$ws.Range("L132").Value = 10224.9375
$ws.Range("M132").Value = -2568.3716
$ws.Range("N132").Value = -15284.9375
$ws.Range("H136").Value = 8193.379000000001
$ws.Range("I136").Value = 9178.913
$ws.Range("J136").Value = 4415.5
$ws.Range("K136").Value = 27536.739
$ws.Range("L136").Value = 13246.5
$ws.Range("M136").Value = -24986.739
$ws.Range("N136").Value = -18346.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2497.513
$ws.Range("I94").Value = 1540.862
$ws.Range("K94").Value = 1540.862
$ws.Range("M94").Value = -1089.862
$ws.Range("H103").Value = 32664
$ws.Range("J103").Value = 32664
$ws.Range("L103").Value = 32664
$ws.Range("N103").Value = -35008
$ws.Range("H134").Value = 5520.1055
$ws.Range("I134").Value = 6101
$ws.Range("J134").Value = 3893.6
$ws.Range("K134").Value = 18303
$ws.Range("L134").Value = 11680.8
$ws.Range("M134").Value = -15768
$ws.Range("N134").Value = -16750.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H105").Value = 668900
$ws.Range("I105").Value = 668900
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 668900
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -667153
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 1810.1177
$ws.Range("I132").Value = 1370.5834
$ws.Range("K132").Value = 4111.7502
$ws.Range("M132").Value = -1581.7502
$ws.Range("H134").Value = 9924.833000000001
$ws.Range("I134").Value = 10554.454
$ws.Range("J134").Value = 2999
$ws.Range("K134").Value = 31663.362
$ws.Range("L134").Value = 8997
$ws.Range("M134").Value = -29128.362
$ws.Range("N134").Value = -14067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 66.2
$ws.Range("I26").Value = 57.75
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 173.25
$ws.Range("L26").Value = 300
$ws.Range("M26").Value = 114.75
$ws.Range("N26").Value = -876
$ws.Range("H87").Value = 14486.385
$ws.Range("I87").Value = 8558.833000000001
$ws.Range("J87").Value = 19567.143
$ws.Range("K87").Value = 25676.499
$ws.Range("L87").Value = 58701.429
$ws.Range("M87").Value = -24428.499
$ws.Range("N87").Value = -61197.429
$ws.Range("H90").Value = 14486.385
$ws.Range("I90").Value = 8558.833000000001
$ws.Range("J90").Value = 19567.143
$ws.Range("K90").Value = 77029.497
$ws.Range("L90").Value = 176104.287
$ws.Range("M90").Value = -70789.497
$ws.Range("N90").Value = -188584.287
$ws.Range("H92").Value = 596.6667
$ws.Range("H98").Value = 2356
$ws.Range("I98").Value = 2138.4
$ws.Range("J98").Value = 2900
$ws.Range("K98").Value = 6415.200000000001
$ws.Range("L98").Value = 8700
$ws.Range("M98").Value = -4917.200000000001
$ws.Range("N98").Value = -11696
$ws.Range("H99").Value = 9281.143
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H122").Value = 5089.75
$ws.Range("I122").Value = 1563.3334
$ws.Range("J122").Value = 5795.033
$ws.Range("K122").Value = 14070.0006
$ws.Range("L122").Value = 52155.29700000001
$ws.Range("M122").Value = -11620.0006
$ws.Range("N122").Value = -57055.29700000001
$ws.Range("H132").Value = 66125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 4259.5
$ws.Range("I36").Value = 4500
$ws.Range("J36").Value = 4019
$ws.Range("K36").Value = 4500
$ws.Range("L36").Value = 4019
$ws.Range("M36").Value = -4015
$ws.Range("N36").Value = -4989
$ws.Range("H44").Value = 29984
$ws.Range("I44").Value = 29950
$ws.Range("J44").Value = 29992.5
$ws.Range("K44").Value = 29950
$ws.Range("L44").Value = 29992.5
$ws.Range("M44").Value = -29354
$ws.Range("N44").Value = -31184.5
$ws.Range("H102").Value = 5622.5713
$ws.Range("I102").Value = 5612.769
$ws.Range("J102").Value = 5750
$ws.Range("K102").Value = 5612.769
$ws.Range("L102").Value = 5750
$ws.Range("M102").Value = -3990.769
$ws.Range("N102").Value = -8994

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 37206.46
$ws.Range("I40").Value = 38307
$ws.Range("K40").Value = 38307
$ws.Range("M40").Value = -38171
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H132").Value = 376335.38
$ws.Range("I132").Value = 535597.4
$ws.Range("J132").Value = 4724.0835
$ws.Range("K132").Value = 1606792.2
$ws.Range("L132").Value = 14172.2505
$ws.Range("M132").Value = -1604262.2
$ws.Range("N132").Value = -19232.2505
$ws.Range("H136").Value = 5223.9375
$ws.Range("I136").Value = 3003.5
$ws.Range("J136").Value = 5541.143
$ws.Range("K136").Value = 9010.5
$ws.Range("L136").Value = 16623.429
$ws.Range("M136").Value = -6460.5
$ws.Range("N136").Value = -21723.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 4916666.5
$ws.Range("I9").Value = 4999999.5
$ws.Range("K9").Value = 4999999.5
$ws.Range("M9").Value = -4999859.5
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H136").Value = 600952.4399999999
$ws.Range("I136").Value = 705475.5600000001
$ws.Range("J136").Value = 26075.25
$ws.Range("K136").Value = 2116426.68
$ws.Range("L136").Value = 78225.75
$ws.Range("M136").Value = -2113876.68
$ws.Range("N136").Value = -83325.75
